$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the individual G-column values that changed
$ws.Range("G7").Value = 4980
$ws.Range("G14").Value = 2728
$ws.Range("G19").Value = 2639
$ws.Range("G26").Value = 2851
$ws.Range("G28").Value = 2689
$ws.Range("G29").Value = 2738
$ws.Range("G30").Value = 2718
$ws.Range("G31").Value = 2890
$ws.Range("G33").Value = 3017
$ws.Range("G35").Value = 3442
$ws.Range("G36").Value = 3672
$ws.Range("G37").Value = 3583
$ws.Range("G38").Value = 3548
$ws.Range("G39").Value = 3293

# Add the new week 48 row
$ws.Range("F40").Value = 48
$ws.Range("G40").Value = 3447
$ws.Range("H40").Value = 3012
$ws.Range("I40").Formula = "=G40-H40"

# Update the view: scroll position and selection
$ws.Range("I41").Select()
$excel.ActiveWindow.ScrollRow = 29
